# Update the "fragile states" aggregate rows (Africa / Rest of the world)
# on the Tab06 sheet to reflect the revised set of underlying countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tab06")

# Row 97: Afrique, Etats fragiles
$ws.Range("C97").Value = 73.668772692307698
$ws.Range("D97").Value = 69.921260769230798
$ws.Range("E97").Value = 77.852894230769195
$ws.Range("F97").Value = 0.88409269230769005
$ws.Range("G97").Value = 62.093014615384597
$ws.Range("H97").Value = 55.0898015384616
$ws.Range("I97").Value = 69.698148846153899
$ws.Range("J97").Value = 0.76230192307692002

# Row 98: RDM, Etats fragiles
$ws.Range("C98").Value = 93.047320769230794
$ws.Range("D98").Value = 92.761195384615405
$ws.Range("E98").Value = 93.315020769230799
$ws.Range("F98").Value = 0.99209307692308002
$ws.Range("G98").Value = 83.795082307692297
$ws.Range("H98").Value = 80.6872969230769
$ws.Range("I98").Value = 86.988493846153901
$ws.Range("J98").Value = 0.92114384615385003
